$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert 10 new (blank) rows starting at row 35. This pushes the existing
#    rows 35-91 down to 45-101 (old row 41 -> 51, old row 71 -> 81, etc.)
# ---------------------------------------------------------------------------
$ws.Rows("35:44").Insert()

# ---------------------------------------------------------------------------
# 2) Fill in the 5 new building-block rows (35-39) for "prob 15".
#    The cell-write order below matters: it reproduces the exact order in
#    which the new strings were first introduced into the shared string
#    table in the target workbook.
# ---------------------------------------------------------------------------
$ws.Range("A35").Value = "c0025"
$ws.Range("B35").Value = "한 꼭짓점이 원의 중심이고 나머지 두 꼭짓점이 원 위의 점인 삼각형은 항상 이등변삼각형이라는 사실을 이용해서 세 각의 관계식을 구합니다."
$ws.Range("A36").Value = "c0026"
$ws.Range("B36").Value = "두 변의 길이와 그 끼인각의 크기가 모두 같으면 합동이라는 사실을 이용해서 두 삼각형이 합동이라는것을 알아냅니다."
$ws.Range("C35").Value = '$\theta_{3}=\dfrac{\pi}{2}+\dfrac{\theta_{2}}{2}$'
$ws.Range("C36").Value = '삼각헝 $\mathrm{O}_{1} \mathrm{O}_{2} \mathrm{~B}$ 와 $\mathrm{O}_{2} \mathrm{O}_{1} \mathrm{D}$; '
$ws.Range("A37").Value = "c0027"
$ws.Range("B37").Value = "한 변이 지름인 원에 내접하는 삼각형은 직각삼각형이라는 사실을 이용해서 피타고라스 정리를 적용해서 빗변인 지름의 길이를 알아냅니다."
$ws.Range("C37").Value = '$3k$;'
$ws.Range("A38").Value = "c0028"
$ws.Range("C38").Value = '$\cos \dfrac{\theta_{1}}{2}$'
$ws.Range("B38").Value = "한 변이 지름인 원에 내접하는 삼각형은 직각삼각형이라는 사실을 이용해서 원하는 삼각비의 값을 알아냅니다."
$ws.Range("A39").Value = "c0029"
$ws.Range("B39").Value = "구하고자 하는 변의 길이를 변수로 해서 주어진 두 변의 길이와 한 각에 대해 코사인법칙을 적용시켜 이차방정식을 세웁니다. "
$ws.Range("C39").Value = '$3 x^{2}-16 k x+21 k^{2}=0$;'

# ---------------------------------------------------------------------------
# 3) Fill in the new row 5 (m0005 triple). This row was previously blank, so
#    no row shift is involved here; the values are simply written in place.
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "m0005"
$ws.Range("B5").Value = "지금까지의 구해진 (가), (나), (다)의 결과를 갖고 문제에서 요구하는 값을 계산합니다."
$ws.Range("C5").Value = '$f(p) \times g(p)$;'

# ---------------------------------------------------------------------------
# 4) Widen columns B and C to fit the newly added (longer) content.
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 146.71428571428572
$ws.Columns("C").ColumnWidth = 90.14285714285714

# ---------------------------------------------------------------------------
# 5) Restore the selection shown in the target workbook.
# ---------------------------------------------------------------------------
$ws.Range("C42").Select()
